$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Common static values shared by every data row in this sheet
$mercadoId = 11
$mercado   = "Vega Monumental Concepción"
$region    = "Bíobío"
$codreg    = 8
$catId     = 100112052
$categoria = "Albahaca"
$variedad  = "Sin especificar"
$calidad   = "Primera"
$unidad    = "$/docena de matas"
$origen    = "Región Metropolitana"
$kgUnid    = 6
$clasif    = "Hortaliza"

# Target data for rows 8-15 (Fecha, Volumen, PrecioMin, PrecioMax, PrecioProm, Precio$/Kg)
$rows = @(
    @{ Row=8;  Fecha=44987; Volumen=130; PMin=4500; PMax=5000; PProm=4692; PKg=782 },
    @{ Row=9;  Fecha=44650; Volumen=130; PMin=3000; PMax=3500; PProm=3308; PKg=551 },
    @{ Row=10; Fecha=44876; Volumen=80;  PMin=6500; PMax=7000; PProm=6812; PKg=1135 },
    @{ Row=11; Fecha=44658; Volumen=180; PMin=2500; PMax=3000; PProm=2778; PKg=463 },
    @{ Row=12; Fecha=44671; Volumen=150; PMin=3500; PMax=4000; PProm=3733; PKg=622 },
    @{ Row=13; Fecha=44685; Volumen=150; PMin=3000; PMax=3500; PProm=3267; PKg=544 },
    @{ Row=14; Fecha=44659; Volumen=90;  PMin=2500; PMax=3000; PProm=2722; PKg=454 },
    @{ Row=15; Fecha=44643; Volumen=90;  PMin=2800; PMax=3000; PProm=2911; PKg=485 }
)

foreach ($r in $rows) {
    $i = $r.Row

    $ws.Cells.Item($i, 1).Value  = $mercadoId
    $ws.Cells.Item($i, 2).Value  = $mercado
    $ws.Cells.Item($i, 3).Value  = $region

    $dCell = $ws.Cells.Item($i, 4)
    $dCell.Value = $r.Fecha
    $dCell.NumberFormat = $ws.Cells.Item(7, 4).NumberFormat

    $ws.Cells.Item($i, 5).Value  = $codreg
    $ws.Cells.Item($i, 6).Value  = $catId
    $ws.Cells.Item($i, 7).Value  = $categoria
    $ws.Cells.Item($i, 8).Value  = $variedad
    $ws.Cells.Item($i, 9).Value  = $calidad
    $ws.Cells.Item($i, 10).Value = $r.Volumen
    $ws.Cells.Item($i, 11).Value = $r.PMin
    $ws.Cells.Item($i, 12).Value = $r.PMax
    $ws.Cells.Item($i, 13).Value = $r.PProm
    $ws.Cells.Item($i, 14).Value = $unidad
    $ws.Cells.Item($i, 15).Value = $origen
    $ws.Cells.Item($i, 16).Value = $r.PKg
    $ws.Cells.Item($i, 17).Value = $kgUnid
    $ws.Cells.Item($i, 18).Value = $clasif
}
